$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the summary cell text (A8, merged A8:D9) from "总结：" to "总结：任务完成"
$ws.Range("A8").Value = "总结：任务完成"

# Update the selection to match the edited cell's merged range
$ws.Range("A8:D9").Select()
